$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "компания легенд"
$ws.Range("C4").Value = "00-00000"
$ws.Range("E4").Value = "РФ"
$ws.Range("B4").Value = "не действует"
$ws.Range("D4").Value = 0

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = 37921

$ws.Range("C8").Select()
